$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (new TPM values)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5380440000000001
$ws.Range("H2").Value = 1.614132
$ws.Range("I2").Value = 0.9003438764610565
$ws.Range("J2").Value = 0.9003438764610566
$ws.Range("Q2").Value = 0.003815090656
$ws.Range("R2").Value = 0.034335815904
$ws.Range("S2").Value = 0.9003438764610565
$ws.Range("T2").Value = 0.9003438764610566

# Row 3 updates (new TPM values)
$ws.Range("G3").Value = 0.05955433333333333
$ws.Range("I3").Value = 0.09965612353894335
$ws.Range("J3").Value = 0.09965612353894336
$ws.Range("Q3").Value = 0.0004222799262222221
$ws.Range("S3").Value = 0.09965612353894335
$ws.Range("T3").Value = 0.09965612353894336
